$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.900.41"
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("D3").Value = "3.403.63"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.39"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.54"
$ws.Range("E6").Value = "  +5.17%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.393"
$ws.Range("E11").Value = "  +4.54%  "
$ws.Range("D12").Value = "3.988.39"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").Value = "  +5.95%  "
$ws.Range("D15").Value = "3.402.69"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.46"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").Value = "61.974.06"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("E18").Value = "  +5.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.49"
$ws.Range("E20").Value = "  +5.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.64"
$ws.Range("E21").Value = "  +10.16%  "
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("E24").Value = "  +14.71%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.72"
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.30"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.17"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.435.25"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.43"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.01"
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  +13.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.786"
$ws.Range("E42").Value = "  +4.92%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.23"
$ws.Range("E44").Value = "  +3.11%  "
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.00"
$ws.Range("E46").Value = "  +6.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.66"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.30"
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("D50").Value = "2.368.91"
$ws.Range("E50").Value = "  +8.92%  "
